{"js": "// In the \"COMPETENCES TECHNIQUES\" block the skill lines are reshuffled:\n//   \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\"  -> becomes\n//       \"Visualisation : tableau\", immediately followed by a brand-new\n//       \"MLOps : ...\" paragraph.\n//   \"Visualisation : tableau\"                         -> becomes\n//       \"Maths : algebra, algorithms\".\n//   \"MLOps : hadoop, spark, ...\"                       -> becomes\n//       \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\".\n//   \"Maths : algebra, algorithms\" (the original, trailing line) is removed.\n//\n// Net result order: Langages, Visualisation, MLOps, Autres, Maths, ML/AI,\n// Bases de donn\u00e9es.\n\nconst body = context.document.body;\n\n// Resolve all four target paragraphs up front (while each search string is\n// still unique in the document), then mutate them via the captured\n// references so later text collisions (e.g. two paragraphs both reading\n// \"Bases de donn\u00e9es : ...\") can't cause a re-search to pick the wrong one.\nconst basesSearch = body.search(\"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\", { matchCase: true });\nconst visuSearch = body.search(\"Visualisation : tableau\", { matchCase: true });\nconst mlopsSearch = body.search(\n  \"MLOps : hadoop, spark, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n  { matchCase: true }\n);\nconst mathsSearch = body.search(\"Maths : algebra, algorithms\", { matchCase: true });\nbasesSearch.load(\"items\");\nvisuSearch.load(\"items\");\nmlopsSearch.load(\"items\");\nmathsSearch.load(\"items\");\nawait context.sync();\n\nconst basesParagraph = basesSearch.items[0].paragraphs.getFirst();\nconst visuParagraph = visuSearch.items[0].paragraphs.getFirst();\nconst mlopsParagraph = mlopsSearch.items[0].paragraphs.getFirst();\nconst mathsParagraph = mathsSearch.items[0].paragraphs.getFirst();\nawait context.sync();\n\n// Drop the trailing \"Maths : algebra, algorithms\" paragraph entirely.\nmathsParagraph.delete();\n\n// \"MLOps : ...\" -> \"Bases de donn\u00e9es : ...\".\nmlopsParagraph.insertText(\"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\", Word.InsertLocation.replace);\n\n// \"Visualisation : tableau\" -> \"Maths : algebra, algorithms\".\nvisuParagraph.insertText(\"Maths : algebra, algorithms\", Word.InsertLocation.replace);\n\n// \"Bases de donn\u00e9es : ...\" -> \"Visualisation : tableau\", plus a new\n// \"MLOps : ...\" paragraph right after it.\nbasesParagraph.insertText(\"Visualisation : tableau\", Word.InsertLocation.replace);\nbasesParagraph.insertParagraph(\n  \"MLOps : hadoop, spark, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# In the \"COMPETENCES TECHNIQUES\" block the skill lines are reshuffled:\n#   \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\"  -> becomes\n#       \"Visualisation : tableau\", immediately followed by a brand-new\n#       \"MLOps : ...\" paragraph.\n#   \"Visualisation : tableau\"                         -> becomes\n#       \"Maths : algebra, algorithms\".\n#   \"MLOps : hadoop, spark, ...\"                       -> becomes\n#       \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\".\n#   \"Maths : algebra, algorithms\" (the original, trailing line) is removed.\n#\n# Net result order: Langages, Visualisation, MLOps, Autres, Maths, ML/AI,\n# Bases de donn\u00e9es.\n\n$d = $word.ActiveDocument\n\n# Locate each of the four target paragraphs by its current text (so the\n# script does not depend on hard-coded paragraph numbers), then grab a\n# fresh handle on each via the document's own Paragraphs collection -\n# Find.Execute leaves the matched Range/Paragraph's cached .Text stale, so\n# re-fetching by Index keeps later reads/writes accurate.\n$findBases = $d.Content\n$findBases.Find.Execute(\"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\") | Out-Null\n$pBases = $d.Paragraphs.Item($findBases.Paragraphs.Item(1).Index)\n\n$findVisu = $d.Content\n$findVisu.Find.Execute(\"Visualisation : tableau\") | Out-Null\n$pVisu = $d.Paragraphs.Item($findVisu.Paragraphs.Item(1).Index)\n\n$findMlops = $d.Content\n$findMlops.Find.Execute(\"MLOps : hadoop, spark, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\") | Out-Null\n$pMlops = $d.Paragraphs.Item($findMlops.Paragraphs.Item(1).Index)\n\n$findMaths = $d.Content\n$findMaths.Find.Execute(\"Maths : algebra, algorithms\") | Out-Null\n$pMaths = $d.Paragraphs.Item($findMaths.Paragraphs.Item(1).Index)\n\n# Each $p* handle tracks its own paragraph from here on, so the edits below\n# can run in any order without needing to re-resolve indices.\n\n# Drop the trailing \"Maths : algebra, algorithms\" paragraph entirely.\n$pMaths.Range.Delete()\n\n# \"MLOps : ...\" -> \"Bases de donn\u00e9es : ...\".\n$pMlops.Range.Text = \"Bases de donn\u00e9es : SQL, MongoDB, Neo4j, Redis\"\n\n# \"Visualisation : tableau\" -> \"Maths : algebra, algorithms\".\n$pVisu.Range.Text = \"Maths : algebra, algorithms\"\n\n# \"Bases de donn\u00e9es : ...\" -> \"Visualisation : tableau\", plus a new\n# \"MLOps : ...\" paragraph right after it.\n$pBases.Range.Text = \"Visualisation : tableau\"\n$pBases.Range.InsertParagraphAfter()\n$newMlopsParagraph = $d.Paragraphs.Item($pBases.Index + 1)\n$newMlopsParagraph.Range.Text = \"MLOps : hadoop, spark, Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\"\n"}
